$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The model list changed: a run with a method not present in the accepted
# list was inserted before the existing "17_07_07_23_04_06_1_1_500_1" row
# (which, together with everything below it, shifts down by one row).
$ws.Rows("19").Insert()

# New row 19 data (unstyled / default formatting, like the other "fresh" rows).
$ws.Range("A19").Value2 = "17_07_07_23_44_39_0_1_500_1"
$ws.Range("B19").Value2 = 0.212115440228391
$ws.Range("C19").Value2 = 0.177082696994985
$ws.Range("D19").Value2 = 0.348478583991745
$ws.Range("E19").Value2 = 0.128600159255576
$ws.Range("F19").Value2 = 0.0450454702134789
$ws.Range("G19").Value2 = 0.0896377651990727
$ws.Range("H19").Value2 = 0.102231876751721
$ws.Range("I19").Value2 = 0.0755164339928472
$ws.Range("J19").Formula = "=SUM(B19:I19)"

# J19 carries its own explicit "no fill" style (distinct from the plain
# default style used by B19:I19) - mirror that by touching and clearing the
# fill so it gets its own style record.
$ws.Range("J19").Interior.ColorIndex = 6
$ws.Range("J19").Interior.Pattern = -4142

# Move the selection like the recorded edit did.
$ws.Range("F22").Select() | Out-Null
